$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.733.80'
$ws.Range("E2").Value = '  -1.21%  '

$ws.Range("D3").Value = '3.370.20'
$ws.Range("E3").Value = '  -0.40%  '

$ws.Range("E4").Value = '  -0.10%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '568.99'
$ws.Range("E5").Value = '  -1.02%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '135.63'
$ws.Range("E6").Value = '  -0.58%  '

$ws.Range("E7").Value = '  +0.03%  '

$ws.Range("D8").Value = '3.368.10'
$ws.Range("E8").Value = '  -0.42%  '

$ws.Range("E9").Value = '  -1.07%  '

$ws.Range("E10").Value = '  +1.28%  '

$ws.Range("E11").Value = '  -3.30%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.378'
$ws.Range("E12").Value = '  -2.90%  '

$ws.Range("D13").Value = '3.944.95'
$ws.Range("E13").Value = '  -0.52%  '

$ws.Range("E14").Value = '  -0.72%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '25.98'
$ws.Range("E15").Value = '  +1.04%  '

$ws.Range("D16").Value = '3.371.96'
$ws.Range("E16").Value = '  -0.45%  '

$ws.Range("E17").Value = '  -3.71%  '

$ws.Range("D18").Value = '60.805.72'
$ws.Range("E18").Value = '  -1.21%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.80'
$ws.Range("E19").Value = '  -1.03%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.72'
$ws.Range("E20").Value = '  -3.03%  '

$ws.Range("E21").Value = '  -1.95%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '370.83'
$ws.Range("E22").Value = '  -1.49%  '

$ws.Range("D23").Value = '3.508.31'
$ws.Range("E23").Value = '  -0.59%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.546'
$ws.Range("E24").Value = '  -1.93%  '

$ws.Range("E25").Value = '  +0.03%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '70.64'
$ws.Range("E26").Value = '  -0.79%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000123'
$ws.Range("E27").Value = '  -2.59%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.174'
$ws.Range("E28").Value = '  +7.94%  '

$ws.Range("E29").Value = '  -8.34%  '

$ws.Range("E30").Value = '  -0.06%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.34'
$ws.Range("E31").Value = '  -2.28%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.01'
$ws.Range("E32").Value = '  -2.56%  '

$ws.Range("E33").Value = '  -2.45%  '

$ws.Range("E34").Value = '  -0.09%  '

$ws.Range("E35").Value = '  -0.85%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.09'
$ws.Range("E36").Value = '  -4.09%  '

$ws.Range("E37").Value = '  -1.02%  '

$ws.Range("E38").Value = '  -1.07%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '164.37'
$ws.Range("E39").Value = '  -0.69%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0758'
$ws.Range("E40").Value = '  -2.67%  '

$ws.Range("B41").Value = 'Stacks'
$ws.Range("C41").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.74'
$ws.Range("E41").Value = '  +1.13%  '

$ws.Range("B42").Value = 'FirstDigitalUSD'
$ws.Range("C42").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.00'
$ws.Range("E42").Value = '  -0.07%  '

$ws.Range("E43").Value = '  -1.08%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '41.85'
$ws.Range("E44").Value = '  +1.03%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '24.96'
$ws.Range("E45").Value = '  +0.66%  '

$ws.Range("E46").Value = '  -2.10%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.15'
$ws.Range("E47").Value = '  -6.19%  '

$ws.Range("D48").Value = '2.509.29'
$ws.Range("E48").Value = '  +7.16%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '23.35'
$ws.Range("E49").Value = '  +3.57%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.74'
$ws.Range("E50").Value = '  -1.53%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.39'
$ws.Range("E51").Value = '  +0.58%  '
